$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Monthly Sales")

# --- Header row (row 1): add "MoM Growth (%)" and "Cumulative Sales" columns,
#     matching the existing bold/bordered/centered header style (copied from B1) ---
$ws.Range("B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)
$ws.Range("C1").Value = "MoM Growth (%)"
$ws.Range("D1").Value = "Cumulative Sales"

# --- Cumulative Sales (column D) for every data row ---
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 1061.84
$ws.Range("D5").Value = 11539.42
$ws.Range("D6").Value = 24210.54
$ws.Range("D7").Value = 40322.78
$ws.Range("D8").Value = 72439.86
$ws.Range("D9").Value = 80664.72
$ws.Range("D10").Value = 89419.66

# --- Month-over-month growth % (column C) ---
# Rows 2-4 have no prior non-zero month to compare against, so they stay blank.
$ws.Range("C5").Value = 886.74
$ws.Range("C6").Value = 20.94
$ws.Range("C7").Value = 27.16
$ws.Range("C8").Value = 99.33
$ws.Range("C9").Value = -74.39
$ws.Range("C10").Value = 6.44

# Highlight the one negative MoM growth value (C9) with a light red/pink fill
$ws.Range("C9").Interior.Color = 10066431

# --- Total row (row 11): new totals for MoM growth and cumulative sales,
#     matching the existing bold Total-row style (copied from B11) ---
$ws.Range("B11").Copy()
$ws.Range("C11:D11").PasteSpecial(-4122)
$ws.Range("C11").Value = 161.04
$ws.Range("D11").Value = 89419.66

$excel.CutCopyMode = $false

Write-Host "Monthly Sales sheet updated with MoM Growth and Cumulative Sales columns"
